$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra rows (4 and 5) first so row count matches the new data (A1:T3)
$ws.Rows("4:5").Delete()

# Row 2: Resolving-Mac -> Ccl12 -> Ccr1 -> ECs
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ccr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.23247666666667
$ws.Range("H2").Value = 60.69743
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04191666666666666
$ws.Range("N2").Value = 0.12575
$ws.Range("O2").Value = 0.007985764192544619
$ws.Range("P2").Value = 0.007985764192544619
$ws.Range("Q2").Value = 0.8480779802777777
$ws.Range("R2").Value = 7.6327018225
$ws.Range("S2").Value = 0.007985764192544619
$ws.Range("T2").Value = 0.007985764192544619

# Row 3: Resolving-Mac -> Ccl12 -> Ccr1 -> Resolving-Mac
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ccr1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.23247666666667
$ws.Range("H3").Value = 60.69743
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.207007
$ws.Range("N3").Value = 15.621021
$ws.Range("O3").Value = 0.9920142358074554
$ws.Range("P3").Value = 0.9920142358074554
$ws.Range("Q3").Value = 105.35064763067
$ws.Range("R3").Value = 948.1558286760298
$ws.Range("S3").Value = 0.9920142358074554
$ws.Range("T3").Value = 0.9920142358074554

Write-Host "done"
